$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value2 = 3867.25
$ws.Range("J42").Value2 = 7397
$ws.Range("L42").Value2 = 22191
$ws.Range("N42").Value2 = -22651
$ws.Range("H86").Value2 = 2362.5
$ws.Range("I86").Value2 = 2125
$ws.Range("K86").Value2 = 2125
$ws.Range("M86").Value2 = -1002
$ws.Range("H89").Value2 = 2362.5
$ws.Range("I89").Value2 = 2125
$ws.Range("K89").Value2 = 10625
$ws.Range("M89").Value2 = -5009
$ws.Range("H98").Value2 = 3267.2856
$ws.Range("J98").Value2 = 4631.4
$ws.Range("L98").Value2 = 4631.4
$ws.Range("N98").Value2 = -7627.4
$ws.Range("H122").Value2 = 3267.2856
$ws.Range("J122").Value2 = 4631.4
$ws.Range("L122").Value2 = 13894.2
$ws.Range("N122").Value2 = -18794.2
$ws.Range("H132").Value2 = 3749.4
$ws.Range("I132").Value2 = 3676.875
$ws.Range("J132").Value2 = 5490
$ws.Range("K132").Value2 = 11030.625
$ws.Range("L132").Value2 = 16470
$ws.Range("M132").Value2 = -8500.625
$ws.Range("N132").Value2 = -21530
$ws.Range("H137").Value2 = 3919.238
$ws.Range("I137").Value2 = 3767.6667
$ws.Range("K137").Value2 = 11303.0001
$ws.Range("M137").Value2 = -8753.000100000001
$ws.Range("H138").Value2 = 5399.7847
$ws.Range("I138").Value2 = 5824
$ws.Range("J138").Value2 = 5364.433
$ws.Range("K138").Value2 = 17472
$ws.Range("L138").Value2 = 16093.299
$ws.Range("M138").Value2 = -12332
$ws.Range("N138").Value2 = -26373.299

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2062.8125
$ws.Range("I2").Value2 = 1949
$ws.Range("K2").Value2 = 1949
$ws.Range("M2").Value2 = -1836
$ws.Range("H45").Value2 = 8452.209000000001
$ws.Range("I45").Value2 = 35991.168
$ws.Range("K45").Value2 = 35991.168
$ws.Range("M45").Value2 = -35614.168
$ws.Range("H50").Value2 = 572.3333
$ws.Range("I50").Value2 = 190.5
$ws.Range("J50").Value2 = 763.25
$ws.Range("K50").Value2 = 190.5
$ws.Range("L50").Value2 = 763.25
$ws.Range("M50").Value2 = 523.5
$ws.Range("N50").Value2 = -2191.25
$ws.Range("H74").Value2 = 5607.3184
$ws.Range("I74").Value2 = 4159.7334
$ws.Range("K74").Value2 = 4159.7334
$ws.Range("M74").Value2 = -3285.7334
$ws.Range("H77").Value2 = 5607.3184
$ws.Range("I77").Value2 = 4159.7334
$ws.Range("K77").Value2 = 20798.667
$ws.Range("M77").Value2 = -16430.667
$ws.Range("H116").Value2 = 2062.8125
$ws.Range("I116").Value2 = 1949
$ws.Range("K116").Value2 = 1949
$ws.Range("M116").Value2 = 345

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2062.8125
$ws.Range("I3").Value2 = 1949
$ws.Range("K3").Value2 = 1949
$ws.Range("M3").Value2 = -1835
$ws.Range("H20").Value2 = 5623
$ws.Range("J20").Value2 = 5553.5
$ws.Range("L20").Value2 = 5553.5
$ws.Range("N20").Value2 = -6047.5
$ws.Range("H86").Value2 = 4042
$ws.Range("I86").Value2 = 3757.4
$ws.Range("J86").Value2 = 4397.75
$ws.Range("K86").Value2 = 3757.4
$ws.Range("L86").Value2 = 4397.75
$ws.Range("M86").Value2 = -2634.4
$ws.Range("N86").Value2 = -6643.75
$ws.Range("H89").Value2 = 4042
$ws.Range("I89").Value2 = 3757.4
$ws.Range("J89").Value2 = 4397.75
$ws.Range("K89").Value2 = 18787
$ws.Range("L89").Value2 = 21988.75
$ws.Range("M89").Value2 = -13171
$ws.Range("N89").Value2 = -33220.75
$ws.Range("H94").Value2 = 916.8261
$ws.Range("I94").Value2 = 1011.6875
$ws.Range("K94").Value2 = 1011.6875
$ws.Range("M94").Value2 = -560.6875
$ws.Range("H99").Value2 = 5799.8
$ws.Range("I99").Value2 = 5777.5557
$ws.Range("K99").Value2 = 5777.5557
$ws.Range("M99").Value2 = -4279.5557
$ws.Range("H105").Value2 = 2182.375
$ws.Range("I105").Value2 = 2291.3333
$ws.Range("K105").Value2 = 2291.3333
$ws.Range("M105").Value2 = -544.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 4144.9473
$ws.Range("I22").Value2 = 4423.7334
$ws.Range("K22").Value2 = 4423.7334
$ws.Range("M22").Value2 = -4073.7334
$ws.Range("H41").Value2 = 32532
$ws.Range("J41").Value2 = 32532
$ws.Range("L41").Value2 = 32532
$ws.Range("N41").Value2 = -33388
$ws.Range("H50").Value2 = 52000
$ws.Range("J50").Value2 = 52000
$ws.Range("L50").Value2 = 52000
$ws.Range("N50").Value2 = -53250
$ws.Range("H58").Value2 = 10188.125
$ws.Range("I58").Value2 = 11009
$ws.Range("K58").Value2 = 11009
$ws.Range("M58").Value2 = -10806
$ws.Range("H134").Value2 = 1837.0555
$ws.Range("I134").Value2 = 1520.6154
$ws.Range("K134").Value2 = 4561.8462
$ws.Range("M134").Value2 = -2026.8462
$ws.Range("H136").Value2 = 10188.125
$ws.Range("I136").Value2 = 11009
$ws.Range("K136").Value2 = 33027
$ws.Range("M136").Value2 = -30477

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value2 = 9000
$ws.Range("I100").Value2 = 0
$ws.Range("J100").Value2 = 9000
$ws.Range("K100").Value2 = 0
$ws.Range("L100").Value2 = 27000
$ws.Range("N100").Value2 = -28622
$ws.Range("H113").Value2 = 4175.7856
$ws.Range("I113").Value2 = 2846.75
$ws.Range("J113").Value2 = 4397.2915
$ws.Range("K113").Value2 = 8540.25
$ws.Range("L113").Value2 = 13191.8745
$ws.Range("M113").Value2 = -6370.25
$ws.Range("N113").Value2 = -17531.8745
$ws.Range("H134").Value2 = 16035.889
$ws.Range("I134").Value2 = 16165.375
$ws.Range("K134").Value2 = 48496.125
$ws.Range("M134").Value2 = -43426.125
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("H76").Value2 = 0
$ws.Range("J76").Value2 = 0
$ws.Range("L76").Value2 = 0
$ws.Range("H79").Value2 = 0
$ws.Range("J79").Value2 = 0
$ws.Range("L79").Value2 = 0
$ws.Range("H97").Value2 = 1399.4
$ws.Range("I97").Value2 = 1399.4
$ws.Range("K97").Value2 = 1399.4
$ws.Range("M97").Value2 = -903.4000000000001
$ws.Range("H132").Value2 = 7457.7173
$ws.Range("I132").Value2 = 7067.706
$ws.Range("J132").Value2 = 8562.75
$ws.Range("K132").Value2 = 21203.118
$ws.Range("L132").Value2 = 25688.25
$ws.Range("M132").Value2 = -18673.118
$ws.Range("N132").Value2 = -30748.25
$ws.Range("N34").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1401.9333
$ws.Range("I16").Value2 = 1444.4166
$ws.Range("K16").Value2 = 1444.4166
$ws.Range("M16").Value2 = -1274.4166
$ws.Range("H122").Value2 = 4982.8335
$ws.Range("I122").Value2 = 4996.143
$ws.Range("J122").Value2 = 4964.2
$ws.Range("K122").Value2 = 14988.429
$ws.Range("L122").Value2 = 14892.6
$ws.Range("M122").Value2 = -12538.429
$ws.Range("N122").Value2 = -19792.6
$ws.Range("H138").Value2 = 76467.875
$ws.Range("I138").Value2 = 12000
$ws.Range("J138").Value2 = 85677.57000000001
$ws.Range("K138").Value2 = 12000
$ws.Range("L138").Value2 = 85677.57000000001
$ws.Range("M138").Value2 = -6860
$ws.Range("N138").Value2 = -95957.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value2 = 155800.67
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 155800.67
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 155800.67
$ws.Range("N5").Value2 = -156024.67
$ws.Range("H132").Value2 = 783.7
$ws.Range("I132").Value2 = 774.25
$ws.Range("K132").Value2 = 2322.75
$ws.Range("M132").Value2 = 207.25
$ws.Range("H136").Value2 = 21470.9
$ws.Range("J136").Value2 = 4977.8
$ws.Range("L136").Value2 = 14933.4
$ws.Range("N136").Value2 = -20033.4
$ws.Range("M5").ClearContents()
